# Naive Bayes Testing implemented
# Adds a new "CompleteNB" results worksheet (mirroring the existing
# CompleteSVM* / CompleteRFC2 results sheets) as the last sheet in the
# workbook, containing the f1_score / accuracy test results.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end
# of the tab order (sheet7 / CompleteNB).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CompleteNB"

# Header row
$ws.Range("B1").Value = "f1_score"
$ws.Range("C1").Value = "accuracy"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

# Match the look of the analogous cells on the other results sheets:
# bold font, centered horizontally, top-aligned vertically, thin box border.
$headerRange = $ws.Range("B1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

$labelCell = $ws.Range("A2")
$labelCell.Font.Bold = $true
$labelCell.HorizontalAlignment = -4108     # xlCenter
$labelCell.VerticalAlignment = -4160       # xlTop
$labelCell.Borders.LineStyle = 1           # xlContinuous
